$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 16
$ws.Range("I16").NumberFormat = "@"
$ws.Range("D16").Value = "image_20250807111344_ppp0.jpg"
$ws.Range("I16").Value = "642,530,686,576"

# Update row 17
$ws.Range("I17").NumberFormat = "@"
$ws.Range("J17").NumberFormat = "@"
$ws.Range("D17").Value = "image_20250807111344_ppp0.jpg"
$ws.Range("I17").Value = "794,481,831,526"
$ws.Range("J17").Value = "0.71"

# Add new row 22
$ws.Range("I22").NumberFormat = "@"
$ws.Range("J22").NumberFormat = "@"
$ws.Range("A22").Value = "66efa766-1456-4beb-b92a-0615a2fc41bb"
$ws.Range("B22").Value = "mosca"
$ws.Range("C22").Value = 45893
$ws.Range("C22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D22").Value = "image_20250824214658_ppp0.jpg"
$ws.Range("E22").Value = "PLACA_20250717165933"
$ws.Range("F22").Value = "Beja"
$ws.Range("G22").Value = 38.02035
$ws.Range("H22").Value = -7.94715
$ws.Range("I22").Value = "1272,293,1315,331"
$ws.Range("J22").Value = "0.69"

$wb.Save()
